# Update gh-pages generated output (456a3b4):
#  - bump "want to go" counts (column F) on several rows
#  - append "西门" (west gate) to the 郎园station venue address (column D)
# Applies to sheet "展览" (exhibitions), sheet "演出" (performances) and the
# aggregated "全部类型" (all types) sheet, which caches the same rows.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

$venueNew = "半截塔路53号首创郎园station西门 郎园station中央车站文化中心"

# ---------------------------------------------------------------------------
# Sheet "展览" (展览 / exhibitions)
# ---------------------------------------------------------------------------
$wsExpo.Range("F4").Value  = 1292
$wsExpo.Range("F6").Value  = 324
$wsExpo.Range("F9").Value  = 7002
$wsExpo.Range("F13").Value = 7896
$wsExpo.Range("F16").Value = 5468

$wsExpo.Range("D17").Value = $venueNew
$wsExpo.Range("D18").Value = $venueNew
$wsExpo.Range("F18").Value = 2354
$wsExpo.Range("F19").Value = 1002
$wsExpo.Range("F21").Value = 283
$wsExpo.Range("F22").Value = 376
$wsExpo.Range("D23").Value = $venueNew

$wsExpo.Range("F26").Value = 238
$wsExpo.Range("F28").Value = 2184
$wsExpo.Range("F30").Value = 251
$wsExpo.Range("F32").Value = 84
$wsExpo.Range("F36").Value = 1442
$wsExpo.Range("F39").Value = 2201

# ---------------------------------------------------------------------------
# Sheet "演出" (演出 / performances)
# ---------------------------------------------------------------------------
$wsShow.Range("F3").Value = 67
$wsShow.Range("F4").Value = 48

# ---------------------------------------------------------------------------
# Sheet "全部类型" (aggregated view of all the other sheets)
# ---------------------------------------------------------------------------
$wsAll.Range("F6").Value  = 1292
$wsAll.Range("F9").Value  = 324
$wsAll.Range("F12").Value = 7002
$wsAll.Range("F16").Value = 7896
$wsAll.Range("F19").Value = 5468

$wsAll.Range("D20").Value = $venueNew
$wsAll.Range("D21").Value = $venueNew
$wsAll.Range("F21").Value = 2354
$wsAll.Range("F22").Value = 1002
$wsAll.Range("F24").Value = 283
$wsAll.Range("F25").Value = 376
$wsAll.Range("D26").Value = $venueNew

$wsAll.Range("F27").Value = 67
$wsAll.Range("F29").Value = 48
$wsAll.Range("F31").Value = 238
$wsAll.Range("F33").Value = 2184
$wsAll.Range("F35").Value = 251
$wsAll.Range("F37").Value = 84
$wsAll.Range("F42").Value = 1442
$wsAll.Range("F45").Value = 2201
